# PackageRenamer_BugDatabase.xlsx edit
#
# Summary of the change:
#   - A new column is inserted before column A (shifting the existing
#     A:E data to B:F).
#   - The new column A becomes a narrow "#" index column: its header
#     cell (row 2) gets the text "#" formatted like the other header
#     cells, and the rest of the column (rows 3-15) is formatted like
#     the rest of the data rows. Column A is given a narrow width.
#   - The title merge (row 1) now spans B1:F1 instead of A1:E1.
#   - The active selection moves to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts all existing
# columns (A:E) one to the right (B:F), carries the values/styles/
# merged cell with them, and keeps the title row merged across the
# same (now-shifted) range.
$ws.Range("A1").EntireColumn.Insert()

# New column A, row 2 is the header for the new "#" column. Copy the
# formatting used by the other header cells (border + centered text)
# from the neighboring header cell, then set its text.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2").Value = "#"

# New column A, rows 3-15 are plain bordered data cells, matching the
# rest of the table body. Copy that formatting from the neighboring
# body cell.
$ws.Range("B3").Copy()
$ws.Range("A3:A15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Make the new index column narrow.
$ws.Columns("A").ColumnWidth = 3.5

# Match the updated selection saved in the workbook.
$ws.Range("B7").Select() | Out-Null
